$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $text) {
    $ws.Range($cell).Value = "'" + $text
    $ws.Range($cell).Style = "Normal"
}

Set-TextCell "D2" "27.161.27"
Set-TextCell "E2" "  +0.91%  "
Set-TextCell "D3" "1.832.20"
Set-TextCell "E3" "  +0.81%  "
Set-TextCell "E4" "  +0.86%  "
Set-TextCell "D5" "312.95"
Set-TextCell "E5" "  +0.91%  "
Set-TextCell "D6" "1.008"
Set-TextCell "E6" "  +0.75%  "
Set-TextCell "D7" "0.4704"
Set-TextCell "E7" "  +1.19%  "
Set-TextCell "D8" "0.3691"
Set-TextCell "E8" "  -0.19%  "
Set-TextCell "D9" "0.07407"
Set-TextCell "E9" "  +0.41%  "
Set-TextCell "D10" "0.8814"
Set-TextCell "E10" "  +1.15%  "
Set-TextCell "D11" "20.43"
Set-TextCell "E11" "  -0.12%  "
Set-TextCell "D12" "1.837.54"
Set-TextCell "E12" "  +0.93%  "
Set-TextCell "D13" "0.07338"
Set-TextCell "E13" "  +3.80%  "
Set-TextCell "D14" "5.473"
Set-TextCell "E14" "  +1.96%  "
Set-TextCell "D15" "92.81"
Set-TextCell "E15" "  +0.59%  "
Set-TextCell "D16" "6.552"
Set-TextCell "E16" "  +0.58%  "
Set-TextCell "E17" "  +0.70%  "
Set-TextCell "D18" "0.000008784"
Set-TextCell "E18" "  +0.69%  "
Set-TextCell "D19" "1.008"
Set-TextCell "E19" "  +0.79%  "
Set-TextCell "D20" "14.79"
Set-TextCell "E20" "  +0.46%  "
Set-TextCell "D21" "27.185.24"
Set-TextCell "E21" "  +0.89%  "
Set-TextCell "D22" "5.306"
Set-TextCell "E22" "  -0.70%  "
Set-TextCell "D23" "10.68"
Set-TextCell "E23" "  +1.29%  "
Set-TextCell "D24" "2.060.59"
Set-TextCell "E24" "  -0.30%  "
Set-TextCell "E25" "  +0.24%  "
Set-TextCell "D26" "152.19"
Set-TextCell "D27" "18.56"
Set-TextCell "E27" "  +0.66%  "
Set-TextCell "D28" "2.156"
Set-TextCell "E28" "  -0.93%  "
Set-TextCell "D29" "5.271"
Set-TextCell "E29" "  -1.06%  "
Set-TextCell "D30" "117.44"
Set-TextCell "E30" "  +1.60%  "
Set-TextCell "D31" "0.08927"
Set-TextCell "E31" "  +0.26%  "
Set-TextCell "D32" "0.7604"
Set-TextCell "E32" "  -0.22%  "
Set-TextCell "D33" "1.171"
Set-TextCell "E33" "  +0.85%  "
Set-TextCell "D34" "4.544"
Set-TextCell "E34" "  +1.20%  "
Set-TextCell "D35" "2.934"
Set-TextCell "E35" "  +0.09%  "
Set-TextCell "D36" "1.009"
Set-TextCell "E36" "  +0.82%  "
Set-TextCell "E37" "  +0.37%  "
Set-TextCell "D38" "0.05336"
Set-TextCell "E38" "  +1.43%  "
Set-TextCell "D39" "0.01961"
Set-TextCell "E39" "  +0.07%  "
Set-TextCell "D40" "2.998"
Set-TextCell "E40" "  +2.32%  "
Set-TextCell "D41" "2.415"
Set-TextCell "E41" "  +1.00%  "
Set-TextCell "D42" "7.321"
Set-TextCell "E42" "  +1.25%  "
Set-TextCell "D43" "0.5345"
Set-TextCell "E43" "  -0.55%  "
Set-TextCell "D44" "0.1663"
Set-TextCell "E44" "  +0.09%  "
Set-TextCell "D45" "8.557"
Set-TextCell "E45" "  +0.61%  "
Set-TextCell "D46" "0.4944"
Set-TextCell "E46" "  -0.23%  "
Set-TextCell "E47" "  +1.35%  "
Set-TextCell "D48" "1.008"
Set-TextCell "E48" "  +0.89%  "
Set-TextCell "D49" "1.670"
Set-TextCell "E49" "  -0.51%  "
Set-TextCell "D50" "103.83"
Set-TextCell "E50" "  +0.77%  "
Set-TextCell "D51" "0.06320"
Set-TextCell "E51" "  +0.53%  "
